$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 93, pushing existing rows 93:212 down to 94:213.
# (Row 212's data therefore lands on row 213, extending the sheet to A1:R213.)
$ws.Rows("93").Insert()

# Populate the newly inserted row 93 with the new observation. Columns
# A,B,C,E,F,G,H,I,N,Q,R hold the same constant values on every data row of
# this sheet, so they are written directly; D,J,K,L,M,O,P carry the new
# per-row figures from the commit.
$ws.Range("A93").Value = 5
$ws.Range("B93").Value = "Macroferia Regional de Talca"
$ws.Range("C93").Value = "Maule"
$ws.Range("D93").Value = 44638
$ws.Range("D93").NumberFormat = $ws.Range("D94").NumberFormat
$ws.Range("E93").Value = 7
$ws.Range("F93").Value = 100112009
$ws.Range("G93").Value = "Acelga"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 200
$ws.Range("K93").Value = 4000
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = 4000
$ws.Range("N93").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O93").Value = "Región del Maule"
$ws.Range("P93").Value = 1000
$ws.Range("Q93").Value = 4
$ws.Range("R93").Value = "Hortaliza"
